# Auto-applies the scheduled market-data refresh to the Leviathan Profits
# workbook: columns H:N on the affected leve rows are recomputed price/profit
# figures. Values (and cell presence - some profit cells are blank when the
# recipe has no HQ/NQ variant) are set to match the refreshed snapshot.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 55: A Real Smooth Move / Lanolin
$ws.Range("H55").Value = 233.2
$ws.Range("I55").Value = 242.5
$ws.Range("J55").Value = 196
$ws.Range("K55").Value = 242.5
$ws.Range("L55").Value = 196
$ws.Range("M55").Value = -28.5
$ws.Range("N55").Value = -624

# Row 57: Quit Your Jib-jab / Gold Needle
$ws.Range("H57").Value = 45895.43
$ws.Range("J57").Value = 45895.43
$ws.Range("L57").Value = 137686.29
$ws.Range("N57").Value = -138684.29

$ws = $wb.Worksheets.Item("ARM")
# Row 5: The Alloyed Truth / Bronze Rivets
$ws.Range("H5").Value = 494.57144
$ws.Range("I5").Value = 494.57144
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 494.57144
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -382.57144

# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 2102.7778
$ws.Range("I61").Value = 1685
$ws.Range("K61").Value = 1685
$ws.Range("M61").Value = -1473

# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 6000
$ws.Range("I122").Value = 6000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 18000
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -15550

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 1580.2222
$ws.Range("J132").Value = 1808.4286
$ws.Range("L132").Value = 5425.2858
$ws.Range("N132").Value = -10485.2858

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 2102.7778
$ws.Range("I136").Value = 1685
$ws.Range("K136").Value = 5055
$ws.Range("M136").Value = -2505

$ws = $wb.Worksheets.Item("BSM")
# Row 4: Mending Fences / Bronze Rivets
$ws.Range("H4").Value = 494.57144
$ws.Range("I4").Value = 494.57144
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 494.57144
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -379.57144

# Row 35: Lancers' Creed / Crowsbeak Hammer
$ws.Range("H35").Value = 34999.715
$ws.Range("J35").Value = 34999.715
$ws.Range("L35").Value = 34999.715
$ws.Range("N35").Value = -35619.715

# Row 54: Get Me to the War on Time / Cobalt Doming Hammer
$ws.Range("H54").Value = 8000
$ws.Range("I54").Value = 8000
$ws.Range("K54").Value = 8000
$ws.Range("M54").Value = -7516

# Row 105: Ingot to Wing It / Molybdenum Ingot
$ws.Range("H105").Value = 9466.4
$ws.Range("I105").Value = 16842.285
$ws.Range("J105").Value = 3012.5
$ws.Range("K105").Value = 16842.285
$ws.Range("L105").Value = 3012.5
$ws.Range("M105").Value = -15095.285
$ws.Range("N105").Value = -6506.5

# Row 107: The Gold Experience / Deepgold Nugget
$ws.Range("H107").Value = 203882.6
$ws.Range("I107").Value = 337466.66
$ws.Range("J107").Value = 3506.5
$ws.Range("K107").Value = 337466.66
$ws.Range("L107").Value = 3506.5
$ws.Range("M107").Value = -335546.66
$ws.Range("N107").Value = -7346.5

$ws = $wb.Worksheets.Item("CRP")
# Row 7: Gridania's Got Talent / Maple Lumber
$ws.Range("H7").Value = 43481024
$ws.Range("I7").Value = 66668764
$ws.Range("J7").Value = 4002.5
$ws.Range("K7").Value = 66668764
$ws.Range("L7").Value = 4002.5
$ws.Range("M7").Value = -66668651
$ws.Range("N7").Value = -4228.5

# Row 16: Raise the Roof / Ash Lumber
$ws.Range("H16").Value = 4330.5713
$ws.Range("I16").Value = 4635.6665
$ws.Range("K16").Value = 4635.6665
$ws.Range("M16").Value = -4348.6665

# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 55576.68
$ws.Range("I31").Value = 65931.875
$ws.Range("J31").Value = 37167.445
$ws.Range("K31").Value = 65931.875
$ws.Range("L31").Value = 37167.445
$ws.Range("M31").Value = -65636.875
$ws.Range("N31").Value = -37757.445

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 55576.68
$ws.Range("I34").Value = 65931.875
$ws.Range("J34").Value = 37167.445
$ws.Range("K34").Value = 65931.875
$ws.Range("L34").Value = 37167.445
$ws.Range("M34").Value = -65729.875
$ws.Range("N34").Value = -37571.445

# Row 52: Spin It Like You Mean It / Mahogany Spinning Wheel
$ws.Range("H52").Value = 51652
$ws.Range("I52").Value = 54967
$ws.Range("J52").Value = 49994.5
$ws.Range("K52").Value = 54967
$ws.Range("L52").Value = 49994.5
$ws.Range("M52").Value = -54673
$ws.Range("N52").Value = -50582.5

# Row 107: Built to Last / White Oak Lumber
$ws.Range("H107").Value = 1536.0294
$ws.Range("I107").Value = 1377.4814
$ws.Range("K107").Value = 1377.4814
$ws.Range("M107").Value = 542.5186000000001

# Row 113: Patient Patients / White Ash Lumber
$ws.Range("H113").Value = 4330.5713
$ws.Range("I113").Value = 4635.6665
$ws.Range("K113").Value = 4635.6665
$ws.Range("M113").Value = -2465.6665

$ws = $wb.Worksheets.Item("CUL")
# Row 64: The Aroma of Faith / Baked Onion Soup
$ws.Range("H64").Value = 3466.3333
$ws.Range("J64").Value = 4000
$ws.Range("L64").Value = 12000
$ws.Range("N64").Value = -12540

# Row 67: Soup's On (L) / Baked Onion Soup
$ws.Range("H67").Value = 3466.3333
$ws.Range("J67").Value = 4000
$ws.Range("L67").Value = 12000
$ws.Range("N67").Value = -13872

# Row 70: Persona non Gratin / Dhalmel Gratin
$ws.Range("H70").Value = 5043.727
$ws.Range("I70").Value = 1661.6666
$ws.Range("J70").Value = 6312
$ws.Range("K70").Value = 4984.9998
$ws.Range("L70").Value = 18936
$ws.Range("M70").Value = -4669.9998
$ws.Range("N70").Value = -19566

# Row 73: Recipe for Disaster (L) / Dhalmel Gratin
$ws.Range("H73").Value = 5043.727
$ws.Range("I73").Value = 1661.6666
$ws.Range("J73").Value = 6312
$ws.Range("K73").Value = 4984.9998
$ws.Range("L73").Value = 18936
$ws.Range("M73").Value = -3892.9998
$ws.Range("N73").Value = -21120

# Row 76: Old Victories, New Tastes / Dhalmel Fricassee
$ws.Range("H76").Value = 57.5
$ws.Range("I76").Value = 57.5
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 172.5
$ws.Range("L76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = 210.5

# Row 79: The Eats of Authenticity (L) / Dhalmel Fricassee
$ws.Range("H79").Value = 57.5
$ws.Range("I79").Value = 57.5
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 172.5
$ws.Range("L79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = 1153.5

# Row 88: Don't Let It Fall Apart / Liver-cheese Sandwich
$ws.Range("H88").Value = 9346.375
$ws.Range("J88").Value = 9346.375
$ws.Range("L88").Value = 28039.125
$ws.Range("N88").Value = -28895.125

# Row 91: Better Come Back with a Sandwich (L) / Liver-cheese Sandwich
$ws.Range("H91").Value = 9346.375
$ws.Range("J91").Value = 9346.375
$ws.Range("L91").Value = 28039.125
$ws.Range("N91").Value = -31003.125

# Row 100: Souper / Gameni
$ws.Range("H100").Value = 8833.833000000001
$ws.Range("I100").Value = 7000
$ws.Range("J100").Value = 9200.6
$ws.Range("K100").Value = 21000
$ws.Range("L100").Value = 27601.8
$ws.Range("M100").Value = -20189
$ws.Range("N100").Value = -29223.8

# Row 105: Fish Box / Chirashi-zushi
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").ClearContents()
$ws.Range("N105").Value = 0

# Row 134: Don't Knock It Till You've Tried It / Mezcal-marinated Swampmonk
$ws.Range("H134").Value = 2457.8572
$ws.Range("I134").Value = 1739.75
$ws.Range("K134").Value = 5219.25
$ws.Range("M134").Value = -149.25

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit / Mythrite Ingot
$ws.Range("H70").Value = 6373.5
$ws.Range("J70").Value = 7250
$ws.Range("L70").Value = 7250
$ws.Range("N70").Value = -7790

# Row 73: Hulls of Broken Dreams (L) / Mythrite Ingot
$ws.Range("H73").Value = 6373.5
$ws.Range("J73").Value = 7250
$ws.Range("L73").Value = 7250
$ws.Range("N73").Value = -9122

# Row 123: Workplace Workout / Ametrine Ring of Fending
$ws.Range("H123").Value = 37265.1
$ws.Range("I123").Value = 17573.25
$ws.Range("J123").Value = 50393
$ws.Range("K123").Value = 17573.25
$ws.Range("L123").Value = 50393
$ws.Range("M123").Value = -15123.25
$ws.Range("N123").Value = -55293

# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 5763
$ws.Range("I126").Value = 7799.2
$ws.Range("K126").Value = 23397.6
$ws.Range("M126").Value = -20927.6

# Row 139: Ringing Gratitude / White Gold Ring of Healing
$ws.Range("H139").Value = 99995
$ws.Range("J139").Value = 99995
$ws.Range("L139").Value = 99995
$ws.Range("N139").Value = -110275

$ws = $wb.Worksheets.Item("LTW")
# Row 55: It's Not a Job, It's a Calling / Peiste Leather
$ws.Range("H55").Value = 329.1579
$ws.Range("I55").Value = 387.66666
$ws.Range("J55").Value = 228.85715
$ws.Range("K55").Value = 387.66666
$ws.Range("L55").Value = 228.85715
$ws.Range("M55").Value = -214.66666
$ws.Range("N55").Value = -574.85715

# Row 61: Spelling Me Softly / Raptor Leather
$ws.Range("H61").Value = 79684.38
$ws.Range("I61").Value = 86166.5
$ws.Range("J61").Value = 1899
$ws.Range("K61").Value = 86166.5
$ws.Range("L61").Value = 1899
$ws.Range("M61").Value = -85964.5
$ws.Range("N61").Value = -2303

# Row 100: Tiger in the Sack / Tiger Leather
$ws.Range("H100").Value = 6912.25
$ws.Range("I100").Value = 6849.6665
$ws.Range("K100").Value = 6849.6665
$ws.Range("M100").Value = -6308.6665

# Row 113: Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 79684.38
$ws.Range("I113").Value = 86166.5
$ws.Range("J113").Value = 1899
$ws.Range("K113").Value = 86166.5
$ws.Range("L113").Value = 1899
$ws.Range("M113").Value = -83996.5
$ws.Range("N113").Value = -6239

# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 109394.69
$ws.Range("I122").Value = 225722.56
$ws.Range("K122").Value = 677167.6799999999
$ws.Range("M122").Value = -674717.6799999999

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 3533.6924
$ws.Range("I132").Value = 3107.8696
$ws.Range("K132").Value = 9323.6088
$ws.Range("M132").Value = -6793.6088

$ws = $wb.Worksheets.Item("WVR")
# Row 26: New Shoes, New Me / Cotton Dress Shoes
$ws.Range("H26").Value = 11749.667
$ws.Range("J26").Value = 15125
$ws.Range("L26").Value = 15125
$ws.Range("N26").Value = -15711

# Row 43: Walk Softly and Carry a Big Halberd / Velveteen Dress Shoes
$ws.Range("H43").Value = 25250
$ws.Range("J43").Value = 25250
$ws.Range("L43").Value = 25250
$ws.Range("N43").Value = -25548

# Row 96: Skills on Display / Ruby Cotton Cloth
$ws.Range("H96").Value = 2619.6
$ws.Range("I96").Value = 2239.8
$ws.Range("K96").Value = 2239.8
$ws.Range("M96").Value = -866.8000000000002

# Row 126: A Polished Purchase / Snow Linen
$ws.Range("H126").Value = 1500
$ws.Range("I126").Value = 1333.3334
$ws.Range("K126").Value = 4000.0002
$ws.Range("M126").Value = -1530.0002

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 9362.823
$ws.Range("I132").Value = 13106.637
$ws.Range("J132").Value = 2499.1667
$ws.Range("K132").Value = 39319.911
$ws.Range("L132").Value = 7497.500100000001
$ws.Range("M132").Value = -36789.911
$ws.Range("N132").Value = -12557.5001
